$d = $word.ActiveDocument

# ===========================================================================
# Change 1 - "...the chain of responsibility pattern provide[s no benefit
# that cannot be achieved more easily through an alternative design patt]
# [_GoBack][ern.]" :
# the stray _GoBack bookmark sitting in the middle of this sentence is
# removed, and the sentence becomes one contiguous run ending in "pattern."
# ===========================================================================
$sentence1 = $d.Content
$sentence1.Find.Execute("s no benefit that cannot be achieved more easily through an alternative design pattern.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Clear then retype the sentence so the (now unbookmarked) neighboring runs
# collapse back together into a single run instead of staying split.
$sentenceRange = $d.Range($sentence1.Start, $sentence1.End)
$sentenceText = $sentenceRange.Text
$sentenceRange.Text = ""
$sentenceRange.InsertAfter($sentenceText)

# ===========================================================================
# Change 2 - last paragraph: "...this pattern is at least as concise as a
# dispatch table, both of which being more concise than a dispatch table."
# becomes "...more concise than a chain of responsibility.", and the
# _GoBack bookmark (tracking the most recent edit) moves to sit right after
# the new final period.
# ===========================================================================
$oldTail = $d.Content
$oldTail.Find.Execute("a dispatch table, both of which being more concise than a dispatch table.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tailStart = $oldTail.Start
$tailEnd = $oldTail.End

$prefix = "a dispatch table, both of which being more concise than a "
$splitPos = $tailStart + $prefix.Length

# Replace "dispatch table." (the final occurrence) with "chain of
# responsibility."
$replaceRange = $d.Range($splitPos, $tailEnd)
$replaceRange.Text = "chain of responsibility."

# The assignment above tends to normalize/merge this run with its
# left-hand neighbor(s). Restore the original two run boundaries --
# between "...concise as " / "a dispatch table...concise than a " / "chain
# of responsibility." -- by briefly bookmarking each boundary point: adding
# a bookmark at a collapsed position splits whichever run currently spans
# it, and that split persists after the temporary bookmark is removed.
foreach ($pos in @($tailStart, $splitPos)) {
    $boundary = $d.Range($pos, $pos)
    $d.Bookmarks.Add("ZZTmpSplit", $boundary)
    $d.Bookmarks("ZZTmpSplit").Delete()
}

# Re-anchor _GoBack at the very end of the document, collapsed, right after
# the new closing period. (A collapsed range that sits exactly on the
# document's last character confuses Bookmarks.Add, so the boundary is
# nudged out of the way with a throw-away character first, then the
# throw-away is removed once the bookmark is safely placed.)
$endPos = $d.Content.End - 1
$nudge = $d.Range($endPos, $endPos)
$nudge.InsertAfter("Z")

$finalSpot = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $finalSpot)

$throwaway = $d.Range($endPos, $endPos + 1)
$throwaway.Text = ""

Write-Output "Applied edits successfully"
